$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 360; this shifts the existing weekly records
# (old rows 360-429) down to rows 362-431, matching the new dimension A1:R431.
$ws.Range("A360:A361").EntireRow.Insert()

# New "Primera" record (row 360) - mirrors the template row (A,B,C,E,F,G,H,I,N,O,Q,R)
# that is now at row 362, with the updated weekly date and price figures.
$ws.Range("A360").Value = 8
$ws.Range("B360").Value = "Terminal La Palmera de La Serena"
$ws.Range("C360").Value = "Coquimbo"
$ws.Range("D360").Value = 45015
$ws.Range("E360").Value = 4
$ws.Range("F360").Value = 100114014
$ws.Range("G360").Value = "Betarraga"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 2100
$ws.Range("K360").Value = 500
$ws.Range("L360").Value = 600
$ws.Range("M360").Value = 550
$ws.Range("N360").Value = "`$/paquete 3 unidades"
$ws.Range("O360").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P360").Value = 183
$ws.Range("Q360").Value = 3
$ws.Range("R360").Value = "Hortaliza"

# New "Segunda" record (row 361)
$ws.Range("A361").Value = 8
$ws.Range("B361").Value = "Terminal La Palmera de La Serena"
$ws.Range("C361").Value = "Coquimbo"
$ws.Range("D361").Value = 45015
$ws.Range("E361").Value = 4
$ws.Range("F361").Value = 100114014
$ws.Range("G361").Value = "Betarraga"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Segunda"
$ws.Range("J361").Value = 1500
$ws.Range("K361").Value = 400
$ws.Range("L361").Value = 450
$ws.Range("M361").Value = 425
$ws.Range("N361").Value = "`$/paquete 3 unidades"
$ws.Range("O361").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P361").Value = 142
$ws.Range("Q361").Value = 3
$ws.Range("R361").Value = "Hortaliza"
